# Product / Sprint backlog update
# Reproduces: insertion of two new "login/logout" stories at the top of the
# Sprint 8 block, re-shuffling of several existing backlog rows, updated
# estimates (2 horas -> 3 horas) and sprint assignment (Sprint 8) for some
# stories, and the trailing "plantilla" row moving from row 43 to row 44.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Historias de Usuario")

# ---------------------------------------------------------------------
# 0) Stash a couple of "format donors" in a scratch area far below the
#    used range so we can re-apply exact existing cell styles later on,
#    regardless of the order in which we overwrite the source rows.
#    (Style indices already present in the sheet: plain bordered cell
#    used by the old blank row 27, and the centered/wrapped style used
#    by the "Estado" column.)
# ---------------------------------------------------------------------
$ws.Range("B29:I29").Copy($ws.Range("B300:I300"))     # style "5" donor (normal data cell)
$ws.Range("B27:I27").Copy($ws.Range("B301:I301"))     # style "11" donor (old blank row27)
$ws.Range("D5").Copy($ws.Range("B302"))                # style "8" donor (centered/wrapped), single cell

# ---------------------------------------------------------------------
# 1) Move the trailing template row from 43 down to 44
# ---------------------------------------------------------------------
$ws.Range("B43:I43").Copy($ws.Range("B44:I44"))
$ws.Range("B43:I43").Clear()

# ---------------------------------------------------------------------
# 2) Shift rows 31-37 down by one (37->38, 36->37, ... 31->32), working
#    from the bottom up so sources are not clobbered before being read.
# ---------------------------------------------------------------------
$ws.Range("B37:I37").Copy($ws.Range("B38:I38"))
$ws.Range("B36:I36").Copy($ws.Range("B37:I37"))
$ws.Range("B35:I35").Copy($ws.Range("B36:I36"))
$ws.Range("B34:I34").Copy($ws.Range("B35:I35"))
$ws.Range("B33:I33").Copy($ws.Range("B34:I34"))
$ws.Range("B32:I32").Copy($ws.Range("B33:I33"))
$ws.Range("B31:I31").Copy($ws.Range("B32:I32"))

# Row 32 (was old row 31, "Notificar inconvenientes en produccion") now
# gets re-estimated from 2 horas to 3 horas.
$ws.Range("F32").Value = "3 horas"

# ---------------------------------------------------------------------
# 3) The old "Documentacion" filler row (old row 28) is relocated to the
#    new row 31, and tagged with Sprint 8.
# ---------------------------------------------------------------------
$ws.Range("B28:I28").Copy($ws.Range("B31:I31"))
$ws.Range("G31").Value = "Sprint 8"

# ---------------------------------------------------------------------
# 4) Rows 29 and 30 stay put but get re-estimated (2 horas -> 3 horas)
#    and assigned to Sprint 8.
# ---------------------------------------------------------------------
$ws.Range("F29").Value = "3 horas"
$ws.Range("G29").Value = "Sprint 8"

$ws.Range("F30").Value = "3 horas"
$ws.Range("G30").Value = "Sprint 8"

# ---------------------------------------------------------------------
# 5) Rows 27 and 28 become two brand new stories (login / logout).
# ---------------------------------------------------------------------

# Row 27: apply the normal data-cell look everywhere, then the
# centered/wrapped "Estado" look on column E only.
$ws.Range("B300:I300").Copy($ws.Range("B27:I27"))
$ws.Range("B302").Copy($ws.Range("E27"))

$ws.Range("B27").Value = ""
$ws.Range("C27").Value = "Como un Administrador, necesito poder iniciar sesion con mis creedenciales"
$ws.Range("D27").Value = "Iniciar sesion "
$ws.Range("E27").Value = "Realizado"
$ws.Range("F27").Value = "2 horas"
$ws.Range("G27").Value = "Sprint 8"
$ws.Range("H27").Value = "Alta"
$ws.Range("I27").Value = ""

# Row 28: columns B, D, I keep the bordered "empty cell" look (style 11,
# borrowed from the scratch donor saved above), column E keeps the
# centered/wrapped "Estado" look, the rest use the normal data style.
$ws.Range("B300:I300").Copy($ws.Range("B28:I28"))
$ws.Range("B301").Copy($ws.Range("B28"))
$ws.Range("B301").Copy($ws.Range("D28"))
$ws.Range("B301").Copy($ws.Range("I28"))
$ws.Range("B302").Copy($ws.Range("E28"))

$ws.Range("B28").Value = ""
$ws.Range("C28").Value = "Como un Administrador, necesito poder cerrar sesion de mi cuenta de usuario"
$ws.Range("D28").Value = "Cerrar sesion"
$ws.Range("E28").Value = "Realizado"
$ws.Range("F28").Value = "2 horas"
$ws.Range("G28").Value = "Sprint 8"
$ws.Range("H28").Value = "Alta"
$ws.Range("I28").Value = ""

# ---------------------------------------------------------------------
# 6) Clean up the scratch/donor area.
# ---------------------------------------------------------------------
$ws.Range("B300:I302").Clear()

# ---------------------------------------------------------------------
# 7) Update the print area to include the new last row (37) and select
#    the same cell/area the author ended up on.
# ---------------------------------------------------------------------
$ws.PageSetup.PrintArea = '$A$1:$I$37'
$ws.Range("C30").Select()
